# Finished Week 13 logging
# Update the Week-13 row (row 2) target-depth data on both the OFF and DEF
# sheets with the latest logged counts.

$wb = $excel.ActiveWorkbook

# --- OFF sheet ---
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 164
$wsOff.Range("C2").Value = 125
$wsOff.Range("D2").Value = 40
$wsOff.Range("E2").Value = 20
$wsOff.Range("G2").Value = 4

# --- DEF sheet ---
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 180
$wsDef.Range("C2").Value = 128
$wsDef.Range("D2").Value = 41
$wsDef.Range("E2").Value = 18
$wsDef.Range("G2").Value = 3
